# "thms results and data added"
#
# The original author cleared out the (numeric "1") values that had been
# populated in the THMS-result columns F:K for most rows, leaving the cell
# formatting/styles in place but with empty cells. Reproduce that by
# clearing the contents (not the formatting) of every one of those cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Every F:K cell (across rows 2-29) that held a "1" value which needs to be
# cleared, while keeping its existing style/formatting intact.
$cellList = "F2,G2,F3,F4,F10,F11,F12,G12,F13,G13,F14,F15,G15,F16,G16,F17,G17," +
            "F18,G18,H18,F19,F20,G20,H20,F21,G21,H21,I21,J21,F22,G22,H22,I22,J22," +
            "F23,G23,H23,I23,J23,F24,G24,H24,F25,G25,H25,I25,J25,F26,G26,H26,I26," +
            "F27,G27,H27,I27,J27,F28,G28,H28,I28,J28,K28,F29,G29,H29,I29,J29"

$targetRange = $ws.Range($cellList)

# ClearContents() on a multi-area Range only touches the first area in this
# host, so walk the Areas collection and clear each contiguous area
# individually -- this removes the <v> while leaving the cell's style (s="")
# attribute untouched, matching the diff exactly.
foreach ($area in $targetRange.Areas) {
    $area.ClearContents()
}
